$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, sorted descending by 1988 value, with the two zero-value rows
# (Russian, Uzbek) dropped entirely.
$data = @(
    @("English", 30.02271828943222),
    @("Japanese", 8.873799625090282),
    @("Spanish", 8.739672369247863),
    @("German", 7.06315094965624),
    @("Chinese", 5.061100643768842),
    @("Arabic", 4.950368506710393),
    @("Portuguese", 4.478074297892415),
    @("French", 4.345483763779167),
    @("Italian", 4.264600137028553),
    @("Malay-Indonesian", 2.152203973861623),
    @("Dutch", 1.847450221762381),
    @("Turkish", 1.413053418210757),
    @("Persian", 1.36833313641462),
    @("Korean", 1.070717162838844),
    @("Polish", 1.027483218452887),
    @("Urdu", 0.7651908511817683),
    @("Thai", 0.7477158898790246),
    @("Swedish", 0.6469766427534291),
    @("Bengali", 0.3993545096946567),
    @("Vietnamese", 0.2672972096240431)
)

# Overwrite the existing rows 2..21 in place (this preserves the formatting
# already applied to column A) with the reordered/trimmed data.
$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# The old data ran through row 23; the new data only needs rows 2..21, so
# remove what's left over (previously Uzbek/Vietnamese, both now obsolete).
$ws.Range("A22:B23").Clear()
